$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Vega Modelo de Temuco - Zanahoria" data block
# (currently occupying rows 479-531), pushing the existing rows down to 481-533.
$ws.Range("A479:R480").EntireRow.Insert()

# --- New row 479 ---
$ws.Cells.Item(479, 1).Value = 10
$ws.Cells.Item(479, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(479, 3).Value = "La Araucanía"
$ws.Cells.Item(479, 4).Value = 45194
$ws.Cells.Item(479, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(479, 5).Value = 9
$ws.Cells.Item(479, 6).Value = 100114013
$ws.Cells.Item(479, 7).Value = "Zanahoria"
$ws.Cells.Item(479, 8).Value = "Sin especificar"
$ws.Cells.Item(479, 9).Value = "Primera"
$ws.Cells.Item(479, 10).Value = 100
$ws.Cells.Item(479, 11).Value = 8000
$ws.Cells.Item(479, 12).Value = 8000
$ws.Cells.Item(479, 13).Value = 8000
$ws.Cells.Item(479, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(479, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(479, 16).Value = 320
$ws.Cells.Item(479, 17).Value = 25
$ws.Cells.Item(479, 18).Value = "Hortaliza"

# --- New row 480 ---
$ws.Cells.Item(480, 1).Value = 10
$ws.Cells.Item(480, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(480, 3).Value = "La Araucanía"
$ws.Cells.Item(480, 4).Value = 45194
$ws.Cells.Item(480, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(480, 5).Value = 9
$ws.Cells.Item(480, 6).Value = 100114013
$ws.Cells.Item(480, 7).Value = "Zanahoria"
$ws.Cells.Item(480, 8).Value = "Sin especificar"
$ws.Cells.Item(480, 9).Value = "Segunda"
$ws.Cells.Item(480, 10).Value = 40
$ws.Cells.Item(480, 11).Value = 6000
$ws.Cells.Item(480, 12).Value = 6000
$ws.Cells.Item(480, 13).Value = 6000
$ws.Cells.Item(480, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(480, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(480, 16).Value = 240
$ws.Cells.Item(480, 17).Value = 25
$ws.Cells.Item(480, 18).Value = "Hortaliza"
